# Generate Report for Handoff
# Adds a new row (for file "6cdf9f30-...") to the Overview, zh-cn and de-de
# sheets of the localization-status workbook, mirroring the existing
# "ddfe76f9-..." row that is already present.

$wb = $excel.ActiveWorkbook

$mdName    = '6cdf9f30-b33d-4a09-9cd6-de8d979a8c49ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdDisplay = 'e2e\6cdf9f30-b33d-4a09-9cd6-de8d979a8c49ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$zhcnXlf   = '6cdf9f30-b33d-4a09-9cd6-de8d979a8c49oooooooooooooooooooooooooooooooooooooooo.f8fcdf8c8f1d39b2f3a3a1c4b8f78cdb165d6a56.zh-cn.xlf'
$dedeXlf   = '6cdf9f30-b33d-4a09-9cd6-de8d979a8c49oooooooooooooooooooooooooooooooooooooooo.f8fcdf8c8f1d39b2f3a3a1c4b8f78cdb165d6a56.de-de.xlf'
$hlUrl     = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1c11ac6ee06c9a878b30a8091873f9a24534f0c/e2e/' + $mdName

$status     = 'Ready for handoff'
$extMd      = '.md'
$emptyStr   = ''
$depFrom    = 'e2e'
$priority   = 'ht'
$falseStr   = 'False'
$trueStr    = 'True'
$epoch      = '0001-01-01 00:00:00'
$hoDate     = '2016-08-30 14:34:32'
$zhHoDate   = '2016-08-30 14:34:28'
$deHoDate   = $hoDate

# ---------------------------------------------------------------
# Overview sheet: add row 3
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hlUrl, "", "", $mdDisplay) | Out-Null
$wsOverview.Range("C3").Value = $extMd
$wsOverview.Range("D3").Value = $emptyStr
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $hoDate

# ---------------------------------------------------------------
# zh-cn sheet: add row 3
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = $extMd
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = $depFrom
$wsZhCn.Range("E3").Value = $priority
$wsZhCn.Range("F3").Value = $falseStr
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $zhHoDate
$wsZhCn.Range("I3").Value = $emptyStr
$wsZhCn.Range("J3").Value = $emptyStr
$wsZhCn.Range("K3").Value = $epoch
$wsZhCn.Range("L3").Value = $emptyStr
$wsZhCn.Range("M3").Value = $trueStr
$wsZhCn.Range("N3").Value = $emptyStr
$wsZhCn.Range("O3").Value = $falseStr
$wsZhCn.Range("P3").Value = $emptyStr

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hlUrl, "", "", $mdName) | Out-Null

# ---------------------------------------------------------------
# de-de sheet: add row 3
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = $extMd
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = $depFrom
$wsDeDe.Range("E3").Value = $priority
$wsDeDe.Range("F3").Value = $falseStr
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $deHoDate
$wsDeDe.Range("I3").Value = $emptyStr
$wsDeDe.Range("J3").Value = $emptyStr
$wsDeDe.Range("K3").Value = $epoch
$wsDeDe.Range("L3").Value = $emptyStr
$wsDeDe.Range("M3").Value = $trueStr
$wsDeDe.Range("N3").Value = $emptyStr
$wsDeDe.Range("O3").Value = $falseStr
$wsDeDe.Range("P3").Value = $emptyStr

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hlUrl, "", "", $mdName) | Out-Null

Write-Host "Applied handoff report rows to Overview, zh-cn, de-de"
